# Applies the "frailty codes.xlsx" update: a new month (Nov-2019, serial 43770)
# of data is added to the "Totals" sheet (row 14), the prior month's row (13)
# is revised with updated figures, a new "all" (API calls, column K) metric is
# added, and the summary rows (30/31) are extended to include the new month.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totals")

# ---------------------------------------------------------------------------
# 1. New shared string / column header "all" in K1
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "all"

# ---------------------------------------------------------------------------
# 2. Row 2: C2 becomes a formula (=+D2) instead of a hard-coded literal, and
#    gains the new K2 (=C2+D2) column.
# ---------------------------------------------------------------------------
$ws.Range("C2").Formula = "=+D2"
$ws.Range("K2").Formula = "=C2+D2"
$ws.Range("K2").NumberFormat = $ws.Range("H2").NumberFormat

# ---------------------------------------------------------------------------
# 3. K3:K14 - new shared "all" column (Not matched + Total matched patients)
# ---------------------------------------------------------------------------
$ws.Range("K3:K14").Formula = "=C3+D3"

# ---------------------------------------------------------------------------
# 4. Row 13 gets revised source figures (B/C/E/F), formulas recalc naturally.
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = 1976
$ws.Range("C13").Value = 61290
$ws.Range("E13").Value = 36794
$ws.Range("F13").Value = 100088

# ---------------------------------------------------------------------------
# 5. Row 14 - brand-new month of data (date already present in A14).
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = 2240
$ws.Range("C14").Value = 61279
$ws.Range("D14").Formula = "=E14+B14"
$ws.Range("E14").Value = 33514
$ws.Range("F14").Value = 97033
$ws.Range("G14").Formula = "=SUM(B`$2:B14)"
$ws.Range("H14").Formula = "=SUM(D`$2:D14)"
$ws.Range("I14").Formula = "=G14/H14"

$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("G14").NumberFormat = $ws.Range("G13").NumberFormat
$ws.Range("H14").NumberFormat = $ws.Range("H13").NumberFormat
$ws.Range("I14").NumberFormat = $ws.Range("I13").NumberFormat
$ws.Range("K14").NumberFormat = $ws.Range("K13").NumberFormat

# ---------------------------------------------------------------------------
# 6. Summary block (rows 30/31): extend ranges to row 14 and clear the stray
#    "100%-B31" formula that used to live in C31.
# ---------------------------------------------------------------------------
$ws.Range("B30").Formula = "=SUM(B2:B14)"
$ws.Range("C30").Formula = "=SUM(D2:D14)"
$ws.Range("D30").Formula = "=SUM(F2:F14)"
$ws.Range("B31").Formula = "=B30/C30"
$ws.Range("C31").ClearContents()

# ---------------------------------------------------------------------------
# 7. Sheet-view bookkeeping: the workbook now opens on the "Totals" tab with
#    C38 selected (instead of the "Frailty graph" chart sheet).
# ---------------------------------------------------------------------------
$ws.Range("C38").Select()
